$wb = $excel.ActiveWorkbook

# --- Sheet management -------------------------------------------------
# Remove the "DeleteList" sheet entirely.
$deleteSheet = $wb.Worksheets.Item("DeleteList")
$deleteSheet.Delete()

# Rename "RegisterList" to "RequestList".
$registerSheet = $wb.Worksheets.Item("RegisterList")
$registerSheet.Name = "RequestList"

# --- Add a header row to the UserList sheet ----------------------------
$ws = $wb.Worksheets.Item("UserList")

# Insert a new blank row at the top; existing rows (and their styles)
# shift down by one.
$ws.Rows.Item(1).Insert()

# Populate the new header row with column titles.
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "First Name"
$ws.Range("C1").Value = "Last Name"
$ws.Range("D1").Value = "DoB"
$ws.Range("E1").Value = "Card #"
$ws.Range("F1").Value = "Last Accessed"
$ws.Range("G1").Value = "Employee Status"
$ws.Range("H1").Value = "Password"

# Update the visible selection to match the saved view.
$ws.Range("B5").Select()
